# Update "Översikt EKERÖ" worksheet:
#  - Column C ("Förändrad") bumps from 46062 -> 46063 for every data row (2..39).
#  - A subset of data rows 2..39 get reshuffled according to a fixed
#    permutation (the underlying case records are the same, they just live
#    at different row positions in the refreshed export). We read every
#    affected source row's current values first (so no Swedish text needs to
#    be hard-coded / re-typed here, avoiding any transcription mistakes),
#    then write the rows back out in their new order, and finally stamp
#    column C on every row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 39

# Columns actually populated on the data rows (others - F, U, Z - are left
# untouched because every row that uses them keeps its original position).
$cols = @(1,2,4,5,7,8,9,10,11,12,13,14,15,16,17,18,19,20,22,23,24,25)

# new row -> old (source) row, only for rows whose content actually moves.
# Rows not listed here (3,4,5,6,11,13,14,16,17,18,19,36,38,2) keep their
# current content untouched aside from the column-C refresh below.
$srcRow = @{
    7=8; 8=9; 9=10; 10=7;
    12=15; 15=12;
    20=23; 21=24; 22=33; 23=25; 24=26; 25=27; 26=32; 27=37; 28=20; 29=21;
    30=39; 31=35; 32=28; 33=30; 34=31; 35=29; 37=34; 39=22
}

# 1) Snapshot the current content of every row that is used as a source,
#    before any writes happen (several rows are both a source and a
#    destination, e.g. the 7/8/9/10 and 12/15 rotations).
$snapshot = @{}
foreach ($source in $srcRow.Values) {
    if (-not $snapshot.ContainsKey($source)) {
        $rowData = @{}
        foreach ($c in $cols) {
            $cell = $ws.Cells.Item($source, $c)
            if ($cell.HasFormula) {
                $rowData[$c] = @{ kind = "formula"; value = $cell.Formula }
            } else {
                $rowData[$c] = @{ kind = "value"; value = $cell.Value2 }
            }
        }
        $snapshot[$source] = $rowData
    }
}

# 2) Write each destination row using the snapshot of its mapped source row.
foreach ($destRow in $srcRow.Keys) {
    $rowData = $snapshot[$srcRow[$destRow]]
    foreach ($c in $cols) {
        $entry = $rowData[$c]
        $cell = $ws.Cells.Item($destRow, $c)
        if ($entry.kind -eq "formula") {
            $cell.Formula = $entry.value
        } else {
            $cell.Value = $entry.value
        }
    }
}

# 3) Column C ("Förändrad") becomes 46063 on every data row.
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 3).Value = 46063
}

# 4) Re-writing the wrapped "Artnamn" column can trigger an auto row-height
#    fit; the source file always pins these rows back to the standard
#    15pt height, so restore it explicitly on every row that was touched
#    (row 39 never had an explicit height to begin with, so it is left
#    alone if untouched).
foreach ($destRow in $srcRow.Keys) {
    if ($destRow -le 38) {
        $ws.Rows.Item($destRow).RowHeight = 15
    }
}
